$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: column C title "Value" -> "Participation" ---
$ws.Cells.Item(1, 3).Value = "Participation"

# --- Column C data rows (1980-2023 in rows 2..45): replace $ amounts with
#     integer participation counts, and restyle to match column B's
#     number format (#,0) instead of the old $ currency format. ---
$participation = @(107,120,93,67,43,60,73,101,95,95,88,88,89,97,93,92,97,100,69,66,68,61,69,50,42,61,51,41,45,42,46,36,36,30,35,29,32,37,23,28,30,22,27,33)

for ($i = 0; $i -lt $participation.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $participation[$i]
    $cell.NumberFormat = "#,0"
}

# --- Total row (row 46): column C total becomes 700 (sum of the new
#     participation counts); both Landings total (col B) and the new
#     Participation total (col C) end up sharing the same bold #,0 style. ---
$totalRow = 46
$ws.Cells.Item($totalRow, 3).Value = 700
$ws.Cells.Item($totalRow, 3).NumberFormat = "#,0"
$ws.Cells.Item($totalRow, 2).NumberFormat = "#,0"
